$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.071.04'
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").Value = '2.512.50'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.83%  '
$ws.Range("E7").Value = '  -0.46%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.562'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.43'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.70'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0817'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("E13").Value = '  +0.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("D15").Value = '2.902.98'
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = '2.518.47'
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.853'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.62%  '
$ws.Range("D18").Value = '47.921.97'
$ws.Range("E18").Value = '  +1.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.48%  '
$ws.Range("D21").Value = '0.0₃0944'
$ws.Range("E21").Value = '  -0.54%  '
$ws.Range("E22").Value = '  +4.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '254.76'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.80%  '
$ws.Range("E25").Value = '  -1.15%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.15%  '
$ws.Range("E29").Value = '  -4.31%  '
$ws.Range("E30").Value = '  +4.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.26'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.08%  '
$ws.Range("E32").Value = '  -0.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.40'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.24%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0787'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.13%  '
$ws.Range("E37").Value = '  -1.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.70'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("E39").Value = '  -1.20%  '
$ws.Range("E40").Value = '  -0.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.29'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '118.71'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.83%  '
$ws.Range("D45").Value = '1.997.82'
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("E46").Value = '  +0.83%  '
$ws.Range("E47").Value = '  -2.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.83'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.10%  '
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '56.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.59%  '
